# Test plan - Tyne & Wear Hospital Finder
# Commit: "Test missed off test plan for when both checkboxes are unchecked. Now added."
#
# 1. On the "Functional tests" sheet (Markers tab test group):
#    - Row 27 (test 18): clarify wording "checked checkbox" -> "checked checkboxes"
#    - Row 28 (test 19): clarify wording "Using the default values" -> "Using the default
#      checked values"; its "Fixed?" column changes from "NR" to "Y" (now considered done)
#    - New row 32 (test 23) added: covers the missing case where both checkboxes are
#      unchecked before clicking Show Markers.
# 2. The active sheet/cell selection moves from the Responsiveness tests sheet to the
#    newly-added row on the Functional tests sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Functional tests")
$ws2 = $wb.Worksheets.Item("Responsiveness tests")

# --- Update wording on existing rows 27 & 28 ---------------------------------
$ws1.Range("B27").Value = "Click the Markers tab to make it active. The tab displays the following: A heading appropriate to the Departure Area, Help text on showing markers, checked checkboxes for both Departure Points and Hospitals, a Show Markers button, and a Clear Map button."

$ws1.Range("B28").Value = "Using the default checked values in both checkboxes, click the Show Markers button. The interactive map displays the Departure Points markers with a red icon, and the Hospital markers with a blue icon."
$ws1.Range("E28").Value = "Y"

# --- Add the new test row (row 32), copying the formatting of row 31 --------
$ws1.Range("A31:E31").Copy()
$ws1.Range("A32:E32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A32").Value = 23
$ws1.Range("B32").Value = "Uncheck both checkboxes and click on Show Markers button. An alert displays warning that at least one checkbox must be checked."
$ws1.Range("C32").Value = "Y"
$ws1.Range("D32").Value = "Nothing happens. No markers are shown and no alert message appears."
$ws1.Range("E32").Value = "Y"

# --- Update the active view: Functional tests sheet / cell E32 --------------
$ws1.Activate()
$ws1.Range("E32").Select()
